# 08/01 - EOD Commit LA
# Adds a new "testingCompanySOI69" variable/value pair and a new
# "idTestingCompanySOI69" variable/value pair to the GeneralVariables sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GeneralVariables")
$ws.Activate()

# Insert a new row above row 5 (after the SOI70 company row) for the
# testingCompanySOI69 / AutoTestingCompany_SOI69 pair.
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).Value = "testingCompanySOI69"
$ws.Cells.Item(5, 2).Value = "AutoTestingCompany_SOI69"

# Insert a new row above what is now row 9 (after the SOI70 id row) for the
# idTestingCompanySOI69 / 0013E00000zXuBpQAK pair.
$ws.Rows.Item(9).Insert()
$ws.Cells.Item(9, 1).Value = "idTestingCompanySOI69"
$ws.Cells.Item(9, 2).Value = "0013E00000zXuBpQAK"

# Update the selection to match the recorded state after the edit.
$ws.Range("A12").Select()
